# Update the "purpose" column (E2:E45) from "S.GISH" to "fullRNASEQ"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2:E45").Value = "fullRNASEQ"

# Update the view: scroll so row 35 is at the top, and select E35:E45 with
# the active cell being E35.
$excel.ActiveWindow.ScrollRow = 35
$ws.Range("E35:E45").Select()
